# Generate Report for Handoff
# Mark the 9da7e663-... file as "Ready for handoff" and refresh its
# Latest Handoff Datetime stamps across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"

# --- Overview sheet ---
# Row 3 corresponds to 9da7e663-14c9-4df4-9d3d-8a93e26a7b00.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value2 = $status
$wsOverview.Range("C3").Value2 = $status
$wsOverview.Range("D3").Value2 = "2016-59-17 02:59:41"

# --- zh-cn sheet ---
# Row 3 corresponds to 9da7e663-14c9-4df4-9d3d-8a93e26a7b00.md
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value2 = $status
$wsZhCn.Range("E3").Value2 = "2016-03-17 02:59:34"

# --- de-de sheet ---
# Row 3 corresponds to 9da7e663-14c9-4df4-9d3d-8a93e26a7b00.md
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value2 = $status
$wsDeDe.Range("E3").Value2 = "2016-03-17 02:59:41"
